$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on 2023-10-07.
# A handful of Price cells are bare decimal numbers (e.g. "211.79"); Excel
# would silently coerce those to Number and drop significant trailing zeros
# (e.g. "65.10" -> 65.1), so those specific cells are pre-formatted as Text
# to preserve the exact source string. Cells containing separators (".") that
# already read back as text, or plainly non-numeric text, are set directly.

$ws.Range("D2").Value = "27.880.73"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.632.43"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.79"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.16"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.866.54"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "1.629.55"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.10"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "27.894.61"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.80"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.30"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.59"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "1.395.45"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").Value = "  +10.76%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.75"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "1.774.99"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.46"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  -3.29%  "

Write-Output "Applied 78 cell updates"
